$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update libraryProtocol (column K) to a uniform value and roboticLibraryPrep
# (column L) to a computed FALSE() formula for every data row (2-27).
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 11).Value = "E7420"
    $ws.Cells.Item($r, 12).Formula = "=FALSE()"
}

# Match the font styling picked up on the updated libraryProtocol column.
$kRange = $ws.Range("K2:K27")
$kRange.Font.Name = "Arial"
$kRange.Font.Size = 11
$kRange.Font.Color = 0

# Leave the selection where it ended up after editing the column.
$ws.Range("K2:K27").Select()
